$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Emp Campus" header (D1) and "BBSR" value (D2) one
# column to the right (E1/E2), preserving the header's bold style.
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("E2").Value = $ws.Range("D2").Value2

# Clear D1's old formatting (it moved to E1) before writing the new header.
$ws.Range("D1").ClearFormats()

# Add the new "Emp Dept" column with its "IT Admin" value in column D.
$ws.Range("D1").Value = "Emp Dept"
$ws.Range("D2").Value = "IT Admin"

# Update selection to E7
$ws.Range("E7").Select()
